$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39; this shifts the existing rows 39-82
# (and all their formatting) down to rows 40-83, exactly like Excel's
# "Insert Row" command.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new weekly record.
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Vega Modelo de Temuco"
$ws.Range("C39").Value = "La Araucanía"
$ws.Range("D39").Value = 44874
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = 100112026
$ws.Range("G39").Value = "Haba"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 35
$ws.Range("K39").Value = 10000
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = 10000
$ws.Range("N39").Value = "$/saco 25 kilos"
$ws.Range("O39").Value = "Región Metropolitana"
$ws.Range("P39").Value = 400
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of the
# column (style index 2 in the original file / "YYYY-MM-DD HH:MM:SS").
$ws.Range("D39").NumberFormat = $ws.Range("D40").NumberFormat
